$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = 45889
$ws.Range("A41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B41").Value = 6395.77978515625
$ws.Range("C41").Value = 6408.39990234375
$ws.Range("D41").Value = 6343.85986328125
$ws.Range("E41").Value = 6406.6201171875
$ws.Range("F41").Value = 4470390000
$ws.Range("G41").Value = -0.0024316693228262
$ws.Range("H41").Value = -1
$ws.Range("I41").Value = 'catchup move lift retail stock a tech loses steam activist investor urge csx to engage in alternative merger discussion dc unemployment rate is the highest in the u for the third straight month patentvest release firstever quantum computing ranking kate spade outlet is selling a shoulder bag for thats a perfect everyday bag intel stock jump softbanks b investment and potential u government stake send share soaring tesla loses bid to kill class action over misleading customer on selfdriving capability for year palantir nvidia and other ai star dim a wall street fall further from it record u appeal court put the brake on contested land transfer for arizona copper mine jerry jones detail how he nearly bought san diego charger using a loan from jimmy hoffa in netflix docuseries meta is shaking up it ai org again techcrunch crypto stock tumble ahead of friday economic conference quartz grocery and drugstore retailer loblaw completes fourforone stock split ap technology summarybrief at pm edt best buy jump on digital marketplace bandwagondoubling it online product may we interest you in a water flight authid auid q earnings call transcript trump weaponization czar urged new york attorney general james to resign over mortgage probe trump weaponization czar urged new york attorney general james to resign over mortgage probe else nutrition schedule first quarter business update conference call bmi settle with radio station to resolve license litigation historic rate increase realtyads expands service portfolio with launch of ainative website service for commercial real estate movie theatre business grew in but still short of prepandemic high statcan dmg blockchain solution inc announces third quarter earnings release date and conference call detail trump administration vying to own a big stake in intel after softbanks billion bet on company trump administration vying to own a big stake in intel after softbanks billion bet on company ap news dmg blockchain solution inc announces third quarter earnings release date and conference call detail blue yonder acquires optoro to streamline return most parent fall into debt for their child survey new andersonville bar clothing boutique open in edgewater and more far north side business news adam ferrari announces phoenix energy scholarship winner the war summit drag on but stock keep rising intels billion rally sends valuation to dotcom level bloombergcom unpaid work for flight attendant could be on the way out a union declares win is strategy becoming a risk to the bitcoin bull market wall street fall on jackson hole jitter reuters calidar announces groundbreaking d mammography system firstinhuman study initiated openai employee stock sale would value chatgpt maker at billion trump weaponization czar urged new york attorney general james to resign over mortgage probe trump weaponization czar urged new york attorney general james to resign over mortgage probe sen tracy pennycuick seeking a responsible budget with sound transit funding opinion crypto stock tumble on tuesday a investor go into riskoff mode pixel watch v pixel watch whats new and different att is directing more manager to relocate or face layoff business insider sofi and lightspark team on blockchainpowered remittance speech by vice chair for supervision bowman on embracing innovation federal reserve board gov gold and mining stock face style resistance signal point to downside risk the new nissan leaf is still one of the cheapest ev in america the verge front porch and ccl hospitality group bring innovative vision for dining to life ap business summarybrief at pm edt ap business summarybrief at pm edt ap business summarybrief at pm edt swatch apologizes for ad showing gesture seen a derogatory and racist pb van leeuwen ice cream coming to westwood village cima and keywest project forge strategic partnership to advance innovation in canada energy sector amer sport raise annual outlook analyst probe china strength hundred of item just got a lot more expensive to import into the u because of trump tariff cnn assessing rocket company insight from financial analyst ap technology summarybrief at pm edt ap technology summarybrief at pm edt ap business summarybrief at pm edt exchangeright secures lease renewal and sustains rent collection for it essential income reit deep dive into helmerich payne stock analyst perspective rating a closer look at mercadolibres option market dynamic spotlight on salesforce analyzing the surge in option activity on semiconductor option trading a deep dive into market sentiment forecasting the future analyst projection for american home rent trump administration vying to own a big stake in intel after softbanks billion bet on company trump administration vying to own a big stake in intel after softbanks billion bet on company trump administration vying to own a big stake in intel after softbanks billion bet on company trump administration vying to own a big stake in intel after softbanks billion bet on company peek inside the plane running the longest flight in the world johnson johnson jnj share purchased by legacy financial advisor inc is the market bullish or bearish on greenwave technology price over earnings overview builder firstsource boulder highway revamp project celebrates progress during st year of work in henderson what doe the market think about aci worldwide looking into applovins recent short interest inszone insurance service expands benefit offering and enters montana market with acquisition of rocky mountain insurance group llc hydro ottawa celebrates year of partnering with capital fair to host special need day tofutti partner with south jersey green and grain to launch plantbased tzatziki sauce despite doge rally and trump optimism blockdags latest miner demo capture market spotlight asia turn to ammonia for power generation union break ground on new town square xrp news today xrp under pressure at may drop toward if support fails cold wallet hype sol eth lead s best crypto project with record presales defi strength shiba inu bull eye gain in august a shib price target rival little pepe lilpepe set for growth athos commerce take a step toward a unified platform athos commerce take a step toward a unified platform over lb of cheese recalled by fda due to contamination concern over serious and sometimes fatal bacteria yahoocom trump expands steel and aluminum tariff to include additional product type cnbc home depot labor day sale beat the rush to grab deal on grill tool appliance and more ford next big ev money move turning battery into business from elevated to worrisome gravity always win billionaire investor say navigating jury selection ethic aba opinion address ai technology jury consultant and client directive walmarts great value shrimp ha risk of radioactive contamination fda say cbs news global tissue engineering and regeneration market set to surge at cagr through stock market today dow give back gain this biotech plunge robinhood launch football betting live coverage investor business daily diddy beat case young thug get his stuff fed target taylor swift ticket more music law dont miss out on aarp perk join now and save on your first year of membership katherine clark back off gaza genocide comment president should be decisive reader commentary lawmaker spends night in texas capitol and could stay longer here why trump weaponization czar urged new york ag to resign over mortgage probe who doe greg abbott think he fooling exclusive medicaid directs state to crack down on illegal immigrant enrollee with monthly check russia open to bilateral or trilateral talk with trump and zelensky putin foreign minister trump doj is investigating whether dc crime stats were manipulated trump weaponization czar urged new york ag james to resign over mortgage probe trump weaponization czar urged new york ag james to resign over mortgage probe trump weaponization czar urged new york ag james to resign over mortgage probe trump weaponization czar urged new york ag james to resign over mortgage probe trump weaponization czar urged new york ag james to resign over mortgage probe trump weaponization czar urged new york ag james to resign over mortgage probe trump weaponization czar urged new york attorney general james to resign over mortgage probe gov josh shapiro speaks on looming septum cut and fare hike video show doj prosecutor ed martin posing for photo outside of new york ag letitia james house boston mayor to trump stop attacking our city to hide your administration failure boston mayor to trump stop attacking our city to hide your administration failure dc unemployment rate is the highest in the u for the third straight month dc unemployment rate is the highest in the u for the third straight month sen tracy pennycuick seeking a responsible budget with sound transit funding opinion costco is at the center of an abortion debate california supreme court save but gut antiarbitration statute the latest trump say no to u troop in ukraine flightsafety international expands arizona presence with new squarefoot learning center in mesa nh judge find once again state is underfunding education takeaway from aps reporting on coach accused of abuse at renowned u gymnastics academy why america is making more artillery shell than ever the fall of florida citrus vance lead effort to block uk from accessing u apple user data resident of a remote island disputed by colombia hope their peruvian government wont forget them port remember when conservative were skeptical of federal overreach ap news summary at pm edt white house trump restored america status a leader of the free world were behind u soldier leading a new catchup crash course say that the army ha to aggressively close the drone gap trump working to arrange summit between putin zelenskyy after monday white house meeting state and dc sue doj to stop immigration requirement on victim fund judge dismisses part of lawsuit over alligator alcatraz immigration detention center judge dismisses part of lawsuit over alligator alcatraz immigration detention center texas gop turn state trooper into armed minder for democrat former u attorney catherine hanaway appointed a missouri next attorney general former u attorney catherine hanaway appointed a missouri next attorney general former u attorney catherine hanaway appointed a missouri next attorney general california republican say redistricting legislation need more review time ask state supreme court to intervene letter dems allow trump to define them this conversation is being recorded trump hot mic moment is the latest in a long global list this conversation is being recorded trump hot mic moment is the latest in a long global list this conversation is being recorded trump hot mic moment is the latest in a long global list this conversation is being recorded trump hot mic moment is the latest in a long global list this conversation is being recorded trump hot mic moment is the latest in a long global list this conversation is being recorded trump hot mic moment is the latest in a long global list immigration cant be handled humanely reader commentary investor alert investigation of soleno therapeutic inc slno announced by holzer holzer llc north africa weekly august james stavridis way to force putin back to the bargaining table mariachi madness boston mayor wu vow to fight ice and prioritize criminal illegal alien over american armstrong williams on crime mayor scott shoot the messenger staff commentary fifth circuit halt west texas am drag show ban a free speech lawsuit continues new hampton man found guilty on poaching charge democratic hopeful should start campaigning now instead of waiting top oversight democrat doj plan to release epstein file in batch a coverup white house launch official tiktok account hillary clinton supreme court will do to gay marriage what they did to abortion polk board of commissioner approve wideranging initiative at august meeting top oversight democrat doj plan to release epstein file in batch a coverup hillary clinton supreme court will do to gay marriage what they did to abortion polk county housing study highlight need for affordable option year after abuse report excoach at renowned u gymnastics academy is arrested by fbi gavin newsom is talking like trump and it infuriating the right gavin newsom is talking like trump and it infuriating the right this conversation is being recorded trump hot mic moment is the latest in a long global list the essential need for recess in elementary school medx holding inc otc medh med x holding q investor update were at an inflection point judge threatened for ruling against trump garrity opts for attack rather than introduction in initial campaign video qatar urge a gaza ceasefire after positive response from hamas the census why trump want a new one voting right act dying a slow death mexico say there no agreement with dea for new border enforcement collaboration senate probe uncovers allegation of widespread abuse in ice custody russia upgraded nuclear battlecruiser back at sea after nearly three decade authority identify family member found dead in murdersuicide at new hampshire home trooper justified in shooting woman in dover ag say sen john kennedy explains how socialist rose to power in the democrat party a only he can watch ice detains driver that struck woman in lewiston the latest trump say no to u troop in ukraine a refreshingly original entry into the mcu federal housing agency adopts englishonly policy rumored presidential hopeful sign bill giving illegals financial aid for college ticket reseller sued for illegally buying over taylor swift ticket the red state push to join the dc occupation nothing say meritocracy like giving dan bongino a babysitter life liberty and school registration gavin newsom want an election in a hurry on his gerrymander here what ha to happen how the tech titan can save social security murdersuicide probed in death of four family member in madbury home a ukrainian sniper just made the world longest kill shot he survived the crucible how the marine set a jm grad on a career path with rochester va clinic rock star give brutal unexpected take on trumpzelenskyy meeting'
